$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTO")

# 1. Remove the cell comment on B9 (and its associated VML drawing / font)
$ws.Range("B9").Comment.Delete()

# 2. Rewrite the transaction rows.
#    Old data occupied rows 2-9 (4 debit/credit pairs); new data only has
#    2 debit/credit pairs (rows 2-5), followed by two blank date rows (6-7).

# Row 2: BS_CREDIT__TRADERECEIVABLECREDITS / Id 1 / 2022-12-31 / UnitId A / Debit / EUR / 100
$ws.Range("A2").Value = "BS_CREDIT__TRADERECEIVABLECREDITS"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "2022-12-31"
$ws.Range("D2").ClearContents()
$ws.Range("F2").Value = "[]"
$ws.Range("G2").Value = "[]"
$ws.Range("H2").Value = "[]"
$ws.Range("I2").Value = "A"
$ws.Range("J2").Value = "BalanceSheet_Debit"
$ws.Range("K2").Value = "EUR"
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = $true
$ws.Range("O2").Value = 1
$ws.Range("P2").ClearContents()
$ws.Range("Q2").Value = 1
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").Value = 100
$ws.Range("U2").Value = "[]"
$ws.Range("V2").Value = "[]"
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()

# Row 3: BS_CASH__BANKACCOUNT_FINANCIALACCOUNT / Id 2 / 2022-12-31 / UnitId A / Credit / EUR / 100
$ws.Range("A3").Value = "BS_CASH__BANKACCOUNT_FINANCIALACCOUNT"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "2022-12-31"
$ws.Range("D3").ClearContents()
$ws.Range("F3").Value = "[]"
$ws.Range("G3").Value = "[]"
$ws.Range("H3").Value = "[]"
$ws.Range("I3").Value = "A"
$ws.Range("J3").Value = "BalanceSheet_Credit"
$ws.Range("K3").Value = "EUR"
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = $true
$ws.Range("O3").Value = 2
$ws.Range("P3").ClearContents()
$ws.Range("Q3").Value = 1
$ws.Range("R3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("T3").Value = 100
$ws.Range("U3").Value = "[]"
$ws.Range("V3").Value = "[]"
$ws.Range("W3").ClearContents()
$ws.Range("X3").ClearContents()

# Row 4: BS_CREDIT__TRADERECEIVABLECREDITS / Id 3 / 2022-12-31 / UnitId A / Debit / EUR / 150
$ws.Range("A4").Value = "BS_CREDIT__TRADERECEIVABLECREDITS"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "2022-12-31"
$ws.Range("D4").ClearContents()
$ws.Range("F4").Value = "[]"
$ws.Range("G4").Value = "[]"
$ws.Range("H4").Value = "[]"
$ws.Range("I4").Value = "A"
$ws.Range("J4").Value = "BalanceSheet_Debit"
$ws.Range("K4").Value = "EUR"
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = $true
$ws.Range("O4").Value = 3
$ws.Range("P4").ClearContents()
$ws.Range("Q4").Value = 2
$ws.Range("R4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").Value = 150
$ws.Range("U4").Value = "[]"
$ws.Range("V4").Value = "[]"
$ws.Range("W4").ClearContents()
$ws.Range("X4").ClearContents()

# Row 5: BS_CASH__BANKACCOUNT_FINANCIALACCOUNT / Id 2 / 2022-12-31 / UnitId A / Credit / EUR / 250 / 150
$ws.Range("A5").Value = "BS_CASH__BANKACCOUNT_FINANCIALACCOUNT"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "2022-12-31"
$ws.Range("D5").ClearContents()
$ws.Range("F5").Value = "[]"
$ws.Range("G5").Value = "[]"
$ws.Range("H5").Value = "[]"
$ws.Range("I5").Value = "A"
$ws.Range("J5").Value = "BalanceSheet_Credit"
$ws.Range("K5").Value = "EUR"
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = $true
$ws.Range("O5").Value = 4
$ws.Range("P5").ClearContents()
$ws.Range("Q5").Value = 2
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("T5").Value = 250
$ws.Range("U5").Value = "[]"
$ws.Range("V5").Value = "[]"
$ws.Range("W5").ClearContents()
$ws.Range("X5").ClearContents()

# Row 6 used to hold a 3rd transaction pair; now blank except the date-formatted cell.
$ws.Range("A6:X6").ClearContents()
$ws.Range("A7:X7").ClearContents()

# Remove the old rows 8-11 entirely so the sheet only spans down to row 7.
$ws.Rows("8:11").Delete()

# Keep the sheet's original column extent (A:AK) even though data only goes to X.
$ws.Range("AK7").Value = "tmp"
$ws.Range("AK7").ClearContents()

# Restore the selected cell to match the saved state in the target workbook.
$ws.Range("A4").Select()
